$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new worksheet "2022-Q1" right after "2021-Q4" (before "总计")
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$newSheet.Name = "2022-Q1"

$template = $wb.Worksheets.Item("2021-Q4")

# Match the outline settings (sheetPr/outlinePr) used by the other quarterly sheets
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

# Match the page setup used by the other quarterly sheets
$newSheet.PageSetup.LeftMargin = $template.PageSetup.LeftMargin
$newSheet.PageSetup.RightMargin = $template.PageSetup.RightMargin
$newSheet.PageSetup.TopMargin = $template.PageSetup.TopMargin
$newSheet.PageSetup.BottomMargin = $template.PageSetup.BottomMargin
$newSheet.PageSetup.HeaderMargin = $template.PageSetup.HeaderMargin
$newSheet.PageSetup.FooterMargin = $template.PageSetup.FooterMargin

# Header row (row 1) - style matches the bold/bordered header used on the other sheets
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "000927"
$newSheet.Range("C2").Value = "博时大中华亚太精选股票(QDII) - 美元现汇"
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "0.32"
$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "92.94"
$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "10.41"
$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.0333"
$newSheet.Range("H2").Value = 1

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").NumberFormat = "@"
$newSheet.Range("B3").Value = "050015"
$newSheet.Range("C3").Value = "博时大中华亚太精选股票(QDII) -人民币"
$newSheet.Range("D3").NumberFormat = "@"
$newSheet.Range("D3").Value = "0.32"
$newSheet.Range("E3").NumberFormat = "@"
$newSheet.Range("E3").Value = "92.94"
$newSheet.Range("F3").NumberFormat = "@"
$newSheet.Range("F3").Value = "10.41"
$newSheet.Range("G3").NumberFormat = "@"
$newSheet.Range("G3").Value = "0.0333"
$newSheet.Range("H3").Value = 1

# Apply the same visual style used by the sibling quarter sheets:
#  - header row (B1:H1) and index column (A2:A3) use the bold bordered style
#  - all other data cells keep the default (unstyled) look
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$template.Range("A2:A3").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)

# Reset number-format-only styling on the text cells back to the plain/default
# style (matches sibling sheets where these cells carry no explicit style),
# while keeping the values stored as text.
$template.Range("C2:C3").Copy()
$newSheet.Range("B2:B3").PasteSpecial(-4122)
$newSheet.Range("D2:G3").PasteSpecial(-4122)
$newSheet.Range("H2:H3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Add the 2022-Q1 summary row to the "总计" sheet (new first data row)
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B3:D3").Copy()
$totalSheet.Range("B2:D2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.07000000000000001

# Re-number the index column (A) for the rows that shifted down
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5

# ---------------------------------------------------------------------------
# 3) Restore the originally active sheet/tab (adding a sheet makes it active)
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
